$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the match scores for row 58 (previously blank), which in turn
# drives the lookup/rank formulas in D58,G58,J58,M58,P58,S58 and the
# season totals in row 70.
$ws.Range("E58").Value = 20
$ws.Range("H58").Value = 100
$ws.Range("K58").Value = 60
$ws.Range("N58").Value = 40
$ws.Range("Q58").Value = 80
$ws.Range("T58").Value = 0

# Update the saved window position to match the author's new view.
$excel.ActiveWindow.Left = 380
$excel.ActiveWindow.Top = 1700

$excel.Calculate()
